# Natmi following Dr Hou advice:
# refresh the Tgm2-Sdc4 LR-pair table (NATMI output) with recomputed
# statistics for the full 6x6 sending-cluster x target-cluster grid
# (ECs, FAPs, M1, M2, Neutro, sCs), growing the sheet from 30 to 36
# data rows (A1:T31 -> A1:T37).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each inner array: Sending, Target, E, F, G, H, I, J, K, L, M, N, O, P, Q, R, S, T
$rows = @(
    @("ECs","ECs","3","1","39.05137633333333","117.154129","0.1500316376309417","0.1500316376309417","2","0.6666666666666666","0.1556706666666667","0.467012","0.00151537522169743","0.00151537522169743","6.079153788060889","54.71238409254801","0.0002273542261366168","0.0002273542261366168"),
    @("ECs","FAPs","3","1","39.05137633333333","117.154129","0.1500316376309417","0.1500316376309417","3","1","37.12743","111.38229","0.3614167567469731","0.3614167567469732","1449.87724121949","13048.89517097541","0.05422394788201208","0.05422394788201208"),
    @("ECs","M1","3","1","39.05137633333333","117.154129","0.1500316376309417","0.1500316376309417","3","1","17.74257933333334","53.227738","0.1727150378838289","0.1727150378838289","692.8721426700224","6235.849284030202","0.02591271997720099","0.02591271997720099"),
    @("ECs","M2","3","1","39.05137633333333","117.154129","0.1500316376309417","0.1500316376309417","3","1","18.87700766666667","56.631023","0.1837581240601467","0.1837581240601467","737.1731304382185","6634.558173943967","0.02756953228073357","0.02756953228073357"),
    @("ECs","Neutro","3","1","39.05137633333333","117.154129","0.1500316376309417","0.1500316376309417","3","1","4.878365","14.635095","0.04748841642225027","0.04748841642225028","190.5068675063616","1714.561807557255","0.007124764884330314","0.007124764884330316"),
    @("ECs","sCs","3","1","39.05137633333333","117.154129","0.1500316376309417","0.1500316376309417","3","1","23.94642","71.83926","0.2331062896651035","0.2331062896651035","935.14065925606","8416.26593330454","0.03497331838052813","0.03497331838052813"),
    @("FAPs","ECs","3","1","3.550472333333333","10.651417","0.01364057374025676","0.01364057374025675","2","0.6666666666666666","0.1556706666666667","0.467012","0.00151537522169743","0.00151537522169743","0.5527043951115557","4.974339556004001","2.067058745572172e-05","2.067058745572172e-05"),
    @("FAPs","FAPs","3","1","3.550472333333333","10.651417","0.01364057374025676","0.01364057374025675","3","1","37.12743","111.38229","0.3614167567469731","0.3614167567469732","131.81991302277","1186.37921720493","0.004929931921371525","0.004929931921371524"),
    @("FAPs","M1","3","1","3.550472333333333","10.651417","0.01364057374025676","0.01364057374025675","3","1","17.74257933333334","53.227738","0.1727150378838289","0.1727150378838289","62.99453704497179","566.950833404746","0.002355932210305608","0.002355932210305607"),
    @("FAPs","M2","3","1","3.550472333333333","10.651417","0.01364057374025676","0.01364057374025675","3","1","18.87700766666667","56.631023","0.1837581240601467","0.1837581240601467","67.02229345662123","603.2006411095911","0.002506566241613681","0.00250656624161368"),
    @("FAPs","Neutro","3","1","3.550472333333333","10.651417","0.01364057374025676","0.01364057374025675","3","1","4.878365","14.635095","0.04748841642225027","0.04748841642225028","17.32049996440167","155.884499679615","0.0006477692460157248","0.0006477692460157247"),
    @("FAPs","sCs","3","1","3.550472333333333","10.651417","0.01364057374025676","0.01364057374025675","3","1","23.94642","71.83926","0.2331062896651035","0.2331062896651035","85.02110169238","765.18991523142","0.003179703533494495","0.003179703533494495"),
    @("M1","ECs","3","1","64.09165766666666","192.274973","0.2462339938068688","0.2462339938068687","2","0.6666666666666666","0.1556706666666667","0.467012","0.00151537522169743","0.00151537522169743","9.977191076741779","89.794719690676","0.0003731368929545273","0.0003731368929545274"),
    @("M1","FAPs","3","1","64.09165766666666","192.274973","0.2462339938068688","0.2462339938068687","3","1","37.12743","111.38229","0.3614167567469731","0.3614167567469732","2379.558533603129","21416.02680242817","0.08899309144253277","0.08899309144253277"),
    @("M1","M1","3","1","64.09165766666666","192.274973","0.2462339938068688","0.2462339938068687","3","1","17.74257933333334","53.227738","0.1727150378838289","0.1727150378838289","1137.151320755675","10234.36188680107","0.04252831356863983","0.04252831356863983"),
    @("M1","M2","3","1","64.09165766666666","192.274973","0.2462339938068688","0.2462339938068687","3","1","18.87700766666667","56.631023","0.1837581240601467","0.1837581240601467","1209.858713143042","10888.72841828738","0.04524749678178799","0.04524749678178799"),
    @("M1","Neutro","3","1","64.09165766666666","192.274973","0.2462339938068688","0.2462339938068687","3","1","4.878365","14.635095","0.04748841642225027","0.04748841642225028","312.6624995530483","2813.962495977435","0.01169326243521438","0.01169326243521438"),
    @("M1","sCs","3","1","64.09165766666666","192.274973","0.2462339938068688","0.2462339938068687","3","1","23.94642","71.83926","0.2331062896651035","0.2331062896651035","1534.76575298222","13812.89177683998","0.05739869268573924","0.05739869268573923"),
    @("M2","ECs","3","1","32.52113333333333","97.5634","0.1249430711472629","0.1249430711472629","2","0.6666666666666666","0.1556706666666667","0.467012","0.00151537522169743","0.00151537522169743","5.062586506755555","45.56327856080001","0.0001893356341393413","0.0001893356341393413"),
    @("M2","FAPs","3","1","32.52113333333333","97.5634","0.1249430711472629","0.1249430711472629","3","1","37.12743","111.38229","0.3614167567469731","0.3614167567469732","1207.426101354","10866.834912186","0.04515651955205008","0.04515651955205008"),
    @("M2","M1","3","1","32.52113333333333","97.5634","0.1249430711472629","0.1249430711472629","3","1","17.74257933333334","53.227738","0.1727150378838289","0.1727150378838289","577.0087881765778","5193.0790935892","0.02157954726652145","0.02157954726652145"),
    @("M2","M2","3","1","32.52113333333333","97.5634","0.1249430711472629","0.1249430711472629","3","1","18.87700766666667","56.631023","0.1837581240601467","0.1837581240601467","613.9016832620222","5525.1151493582","0.02295930436833448","0.02295930436833448"),
    @("M2","Neutro","3","1","32.52113333333333","97.5634","0.1249430711472629","0.1249430711472629","3","1","4.878365","14.635095","0.04748841642225027","0.04748841642225028","158.6499586136667","1427.849627523","0.005933348591716066","0.005933348591716066"),
    @("M2","sCs","3","1","32.52113333333333","97.5634","0.1249430711472629","0.1249430711472629","3","1","23.94642","71.83926","0.2331062896651035","0.2331062896651035","778.7647176759999","7008.882459083999","0.0291250157345015","0.0291250157345015"),
    @("Neutro","ECs","3","1","117.3842696666667","352.152809","0.4509790912329367","0.4509790912329367","2","0.6666666666666666","0.1556706666666667","0.467012","0.00151537522169743","0.00151537522169743","18.27328751518978","164.459587636708","0.000683402540358017","0.000683402540358017"),
    @("Neutro","FAPs","3","1","117.3842696666667","352.152809","0.4509790912329367","0.4509790912329367","3","1","37.12743","111.38229","0.3614167567469731","0.3614167567469732","4358.17625515029","39223.58629635262","0.1629914005141053","0.1629914005141053"),
    @("Neutro","M1","3","1","117.3842696666667","352.152809","0.4509790912329367","0.4509790912329367","3","1","17.74257933333334","53.227738","0.1727150378838289","0.1727150378838289","2082.699717046227","18744.29745341605","0.07789087082711142","0.07789087082711141"),
    @("Neutro","M2","3","1","117.3842696666667","352.152809","0.4509790912329367","0.4509790912329367","3","1","18.87700766666667","56.631023","0.1837581240601467","0.1837581240601467","2215.863758443735","19942.77382599361","0.08287107179531422","0.08287107179531421"),
    @("Neutro","Neutro","3","1","117.3842696666667","352.152809","0.4509790912329367","0.4509790912329367","3","1","4.878365","14.635095","0.04748841642225027","0.04748841642225028","572.6433126924284","5153.789814231855","0.0214162828821977","0.0214162828821977"),
    @("Neutro","sCs","3","1","117.3842696666667","352.152809","0.4509790912329367","0.4509790912329367","3","1","23.94642","71.83926","0.2331062896651035","0.2331062896651035","2810.93302283126","25298.39720548134","0.1051260626738501","0.1051260626738501"),
    @("sCs","ECs","3","1","3.688700333333333","11.066101","0.01417163244173325","0.01417163244173325","2","0.6666666666666666","0.1556706666666667","0.467012","0.00151537522169743","0.00151537522169743","0.5742224400235556","5.168001960212","2.147534065320601e-05","2.147534065320601e-05"),
    @("sCs","FAPs","3","1","3.688700333333333","11.066101","0.01417163244173325","0.01417163244173325","3","1","37.12743","111.38229","0.3614167567469731","0.3614167567469732","136.95196341681","1232.56767075129","0.005121865434901418","0.005121865434901418"),
    @("sCs","M1","3","1","3.688700333333333","11.066101","0.01417163244173325","0.01417163244173325","3","1","17.74257933333334","53.227738","0.1727150378838289","0.1727150378838289","65.44705830105978","589.023524709538","0.002447654034049657","0.002447654034049657"),
    @("sCs","M2","3","1","3.688700333333333","11.066101","0.01417163244173325","0.01417163244173325","3","1","18.87700766666667","56.631023","0.1837581240601467","0.1837581240601467","69.63162447236923","626.6846202513229","0.002604152592362818","0.002604152592362818"),
    @("sCs","Neutro","3","1","3.688700333333333","11.066101","0.01417163244173325","0.01417163244173325","3","1","4.878365","14.635095","0.04748841642225027","0.04748841642225028","17.99482660162166","161.953439414595","0.0006729883827760999","0.0006729883827760999"),
    @("sCs","sCs","3","1","3.688700333333333","11.066101","0.01417163244173325","0.01417163244173325","3","1","23.94642","71.83926","0.2331062896651035","0.2331062896651035","88.33116743613999","794.98050692526","0.003303496656990048","0.003303496656990047")
)

$startRow = 2
for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $startRow + $i
    $data = $rows[$i]
    $ws.Cells.Item($r, 1).Value = $data[0]
    $ws.Cells.Item($r, 2).Value = "Tgm2"
    $ws.Cells.Item($r, 3).Value = "Sdc4"
    $ws.Cells.Item($r, 4).Value = $data[1]
    for ($c = 2; $c -lt $data.Count; $c++) {
        $ws.Cells.Item($r, $c + 3).Value = [double]$data[$c]
    }
}